$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (A and B) as close as the COM ColumnWidth rounding allows
$ws.Range("A1").ColumnWidth = 14.666666666666666
$ws.Range("B1").ColumnWidth = 13.833333333333334

# Update cell values
$ws.Range("A1").Value = -0.0090821767129363965
$ws.Range("B1").Value = 0.0090821764594194825

$ws.Range("A2").Value = -0.0021760635564545558
$ws.Range("B2").Value = 0.0021760632770192754

$ws.Range("A3").Value = 0.023443850904779635
$ws.Range("B3").Value = -0.023443851162496917

$ws.Range("A4").Value = -0.0062379540872743843
$ws.Range("B4").Value = 0.0062379538101951576
